$d = $word.ActiveDocument

# The picture lives in the very last paragraph of the document body
# (immediately before the final sectPr).
$picPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# Remove the centered alignment from that paragraph.
$picPara.Alignment = 0

$picRange = $picPara.Range

# Insert a new paragraph right after the picture paragraph, before touching
# NoProofing, so the new paragraph's run does not inherit <w:noProof/>.
$picRange.InsertParagraphAfter() | Out-Null

# Mark the run that holds the drawing as "no proofing" (<w:noProof/>).
$picRange.NoProofing = 1

# Fill in the newly inserted trailing paragraph with the TP2 update text.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "TP2 Update: I have not made any explicit changes to my design proposal since TP1"
